$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ALC row 64
$ws1.Range("H64").Value = 3777.2273
$ws1.Range("I64").Value = 4124.75
$ws1.Range("J64").Value = 3700
$ws1.Range("K64").Value = 4124.75
$ws1.Range("L64").Value = 3700
$ws1.Range("M64").Value = -3876.75
$ws1.Range("N64").Value = -4196

# ALC row 67
$ws1.Range("H67").Value = 3777.2273
$ws1.Range("I67").Value = 4124.75
$ws1.Range("J67").Value = 3700
$ws1.Range("K67").Value = 4124.75
$ws1.Range("L67").Value = 3700
$ws1.Range("M67").Value = -3266.75
$ws1.Range("N67").Value = -5416

# ALC row 70
$ws1.Range("H70").Value = 16768027
$ws1.Range("I70").Value = 33534314
$ws1.Range("J70").Value = 1740
$ws1.Range("K70").Value = 100602942
$ws1.Range("L70").Value = 5220
$ws1.Range("M70").Value = -100602672
$ws1.Range("N70").Value = -5760

# ALC row 73
$ws1.Range("H73").Value = 16768027
$ws1.Range("I73").Value = 33534314
$ws1.Range("J73").Value = 1740
$ws1.Range("K73").Value = 100602942
$ws1.Range("L73").Value = 5220
$ws1.Range("M73").Value = -100602006
$ws1.Range("N73").Value = -7092

# ALC row 112
$ws1.Range("H112").Value = 2407.3333
$ws1.Range("J112").Value = 2451.3125
$ws1.Range("L112").Value = 7353.9375
$ws1.Range("N112").Value = -9569.9375

# ALC row 129
$ws1.Range("H129").Value = 897.98
$ws1.Range("J129").Value = 897.98
$ws1.Range("L129").Value = 2693.94
$ws1.Range("N129").Value = -12693.94

# ALC row 137
$ws1.Range("H137").Value = 1977
$ws1.Range("I137").Value = 1388.8518
$ws1.Range("J137").Value = 2612.2
$ws1.Range("K137").Value = 4166.555399999999
$ws1.Range("L137").Value = 7836.599999999999
$ws1.Range("M137").Value = -1616.555399999999
$ws1.Range("N137").Value = -12936.6

# ALC row 138
$ws1.Range("H138").Value = 3605.902
$ws1.Range("I138").Value = 2628.7727
$ws1.Range("J138").Value = 3913
$ws1.Range("K138").Value = 7886.3181
$ws1.Range("L138").Value = 11739
$ws1.Range("M138").Value = -2746.3181
$ws1.Range("N138").Value = -22019

# ARM row 32
$ws2.Range("H32").Value = 10087.127
$ws2.Range("I32").Value = 11080.492
$ws2.Range("J32").Value = 5475.0713
$ws2.Range("K32").Value = 11080.492
$ws2.Range("L32").Value = 5475.0713
$ws2.Range("M32").Value = -10793.492
$ws2.Range("N32").Value = -6049.0713

# ARM row 45
$ws2.Range("H45").Value = 1409.375
$ws2.Range("I45").Value = 1375
$ws2.Range("J45").Value = 1650
$ws2.Range("K45").Value = 1375
$ws2.Range("L45").Value = 1650
$ws2.Range("M45").Value = -998
$ws2.Range("N45").Value = -2404

# ARM row 61
$ws2.Range("H61").Value = 2147.2632
$ws2.Range("I61").Value = 2040.64
$ws2.Range("J61").Value = 2352.3076
$ws2.Range("K61").Value = 2040.64
$ws2.Range("L61").Value = 2352.3076
$ws2.Range("M61").Value = -1828.64
$ws2.Range("N61").Value = -2776.3076

# ARM row 74
$ws2.Range("H74").Value = 1273.0975
$ws2.Range("I74").Value = 1092.8485
$ws2.Range("K74").Value = 1092.8485
$ws2.Range("M74").Value = -218.8485000000001

# ARM row 77
$ws2.Range("H77").Value = 1273.0975
$ws2.Range("I77").Value = 1092.8485
$ws2.Range("K77").Value = 5464.2425
$ws2.Range("M77").Value = -1096.2425

# ARM row 122
$ws2.Range("H122").Value = 6511.2256
$ws2.Range("I122").Value = 6593.852
$ws2.Range("J122").Value = 5953.5
$ws2.Range("K122").Value = 19781.556
$ws2.Range("L122").Value = 17860.5
$ws2.Range("M122").Value = -17331.556
$ws2.Range("N122").Value = -22760.5

# ARM row 132
$ws2.Range("H132").Value = 8338.617
$ws2.Range("I132").Value = 8593.056
$ws2.Range("K132").Value = 25779.168
$ws2.Range("M132").Value = -23249.168

# ARM row 136
$ws2.Range("H136").Value = 2147.2632
$ws2.Range("I136").Value = 2040.64
$ws2.Range("J136").Value = 2352.3076
$ws2.Range("K136").Value = 6121.92
$ws2.Range("L136").Value = 7056.9228
$ws2.Range("M136").Value = -3571.92
$ws2.Range("N136").Value = -12156.9228

# BSM row 114
$ws3.Range("H114").Value = 25000
$ws3.Range("J114").Value = 25000
$ws3.Range("L114").Value = 25000
$ws3.Range("N114").Value = -33678

# BSM row 134
$ws3.Range("H134").Value = 2967.0908
$ws3.Range("I134").Value = 2774.8572
$ws3.Range("J134").Value = 3303.5
$ws3.Range("K134").Value = 8324.571599999999
$ws3.Range("L134").Value = 9910.5
$ws3.Range("M134").Value = -5789.571599999999
$ws3.Range("N134").Value = -14980.5

# CRP row 31
$ws4.Range("H31").Value = 2407.8838
$ws4.Range("I31").Value = 2690.3635
$ws4.Range("J31").Value = 2111.9524
$ws4.Range("K31").Value = 2690.3635
$ws4.Range("L31").Value = 2111.9524
$ws4.Range("M31").Value = -2395.3635
$ws4.Range("N31").Value = -2701.9524

# CRP row 34
$ws4.Range("H34").Value = 2407.8838
$ws4.Range("I34").Value = 2690.3635
$ws4.Range("J34").Value = 2111.9524
$ws4.Range("K34").Value = 2690.3635
$ws4.Range("L34").Value = 2111.9524
$ws4.Range("M34").Value = -2488.3635
$ws4.Range("N34").Value = -2515.9524

# CRP row 86
$ws4.Range("H86").Value = 3859.7273
$ws4.Range("I86").Value = 2522.4285
$ws4.Range("J86").Value = 6200
$ws4.Range("K86").Value = 2522.4285
$ws4.Range("L86").Value = 6200
$ws4.Range("M86").Value = -1399.4285
$ws4.Range("N86").Value = -8446

# CRP row 89
$ws4.Range("H89").Value = 3859.7273
$ws4.Range("I89").Value = 2522.4285
$ws4.Range("J89").Value = 6200
$ws4.Range("K89").Value = 12612.1425
$ws4.Range("L89").Value = 31000
$ws4.Range("M89").Value = -6996.1425
$ws4.Range("N89").Value = -42232

# CRP row 122
$ws4.Range("H122").Value = 3591.25
$ws4.Range("I122").Value = 4527.8335
$ws4.Range("J122").Value = 781.5
$ws4.Range("K122").Value = 13583.5005
$ws4.Range("L122").Value = 2344.5
$ws4.Range("M122").Value = -11133.5005
$ws4.Range("N122").Value = -7244.5

# CUL row 68
$ws5.Range("H68").Value = 209239.38
$ws5.Range("I68").Value = 257075.64
$ws5.Range("J68").Value = 1948.8889
$ws5.Range("K68").Value = 771226.92
$ws5.Range("L68").Value = 5846.6667
$ws5.Range("M68").Value = -770415.92
$ws5.Range("N68").Value = -7468.6667

# CUL row 71
$ws5.Range("H71").Value = 209239.38
$ws5.Range("I71").Value = 257075.64
$ws5.Range("J71").Value = 1948.8889
$ws5.Range("K71").Value = 2313680.76
$ws5.Range("L71").Value = 17540.0001
$ws5.Range("M71").Value = -2309624.76
$ws5.Range("N71").Value = -25652.0001

# GSM row 80
$ws6.Range("H80").Value = 4272.143
$ws6.Range("I80").Value = 4381
$ws6.Range("J80").Value = 4000
$ws6.Range("K80").Value = 4381
$ws6.Range("L80").Value = 4000
$ws6.Range("M80").Value = -3383
$ws6.Range("N80").Value = -5996

# GSM row 83
$ws6.Range("H83").Value = 4272.143
$ws6.Range("I83").Value = 4381
$ws6.Range("J83").Value = 4000
$ws6.Range("K83").Value = 21905
$ws6.Range("L83").Value = 20000
$ws6.Range("M83").Value = -16913
$ws6.Range("N83").Value = -29984

# GSM row 102
$ws6.Range("H102").Value = 2099.2195
$ws6.Range("I102").Value = 2387.9312
$ws6.Range("J102").Value = 1401.5
$ws6.Range("K102").Value = 2387.9312
$ws6.Range("L102").Value = 1401.5
$ws6.Range("M102").Value = -765.9312
$ws6.Range("N102").Value = -4645.5

# GSM row 113
$ws6.Range("H113").Value = 2566.2727
$ws6.Range("I113").Value = 1450
$ws6.Range("J113").Value = 2814.3333
$ws6.Range("K113").Value = 1450
$ws6.Range("L113").Value = 2814.3333
$ws6.Range("M113").Value = 720
$ws6.Range("N113").Value = -7154.3333

# LTW row 16
$ws7.Range("H16").Value = 2240.75
$ws7.Range("J16").Value = 3307.8572
$ws7.Range("L16").Value = 3307.8572
$ws7.Range("N16").Value = -3647.8572

# LTW row 40
$ws7.Range("H40").Value = 3213.9285
$ws7.Range("I40").Value = 2700
$ws7.Range("J40").Value = 3499.4443
$ws7.Range("K40").Value = 2700
$ws7.Range("L40").Value = 3499.4443
$ws7.Range("M40").Value = -2564
$ws7.Range("N40").Value = -3771.4443

# LTW row 42
$ws7.Range("H42").Value = 34231.25
$ws7.Range("J42").Value = 29962.5
$ws7.Range("L42").Value = 29962.5
$ws7.Range("N42").Value = -31088.5

# LTW row 49
$ws7.Range("H49").Value = 34231.25
$ws7.Range("J49").Value = 29962.5
$ws7.Range("L49").Value = 29962.5
$ws7.Range("N49").Value = -30256.5

# LTW row 68
$ws7.Range("H68").Value = 2625.2144
$ws7.Range("I68").Value = 1861.1111
$ws7.Range("K68").Value = 1861.1111
$ws7.Range("M68").Value = -1112.1111

# LTW row 71
$ws7.Range("H71").Value = 2625.2144
$ws7.Range("I71").Value = 1861.1111
$ws7.Range("K71").Value = 9305.5555
$ws7.Range("M71").Value = -5561.5555

# LTW row 122
$ws7.Range("H122").Value = 20459986
$ws7.Range("I122").Value = 31254280
$ws7.Range("J122").Value = 14291820
$ws7.Range("K122").Value = 93762840
$ws7.Range("L122").Value = 42875460
$ws7.Range("M122").Value = -93760390
$ws7.Range("N122").Value = -42880360

# WVR row 62
$ws8.Range("H62").Value = 5374.75
$ws8.Range("I62").Value = 5999.6665
$ws8.Range("J62").Value = 4999.8
$ws8.Range("K62").Value = 5999.6665
$ws8.Range("L62").Value = 4999.8
$ws8.Range("M62").Value = -5375.6665
$ws8.Range("N62").Value = -6247.8

# WVR row 65
$ws8.Range("H65").Value = 5374.75
$ws8.Range("I65").Value = 5999.6665
$ws8.Range("J65").Value = 4999.8
$ws8.Range("K65").Value = 29998.3325
$ws8.Range("L65").Value = 24999
$ws8.Range("M65").Value = -26878.3325
$ws8.Range("N65").Value = -31239

# WVR row 81
$ws8.Range("H81").Value = 42340.32
$ws8.Range("I81").Value = 73829
$ws8.Range("J81").Value = 2263.818
$ws8.Range("K81").Value = 147658
$ws8.Range("L81").Value = 4527.636
$ws8.Range("M81").Value = -146597
$ws8.Range("N81").Value = -6649.636

# WVR row 84
$ws8.Range("H84").Value = 42340.32
$ws8.Range("I84").Value = 73829
$ws8.Range("J84").Value = 2263.818
$ws8.Range("K84").Value = 738290
$ws8.Range("L84").Value = 22638.18
$ws8.Range("M84").Value = -732986
$ws8.Range("N84").Value = -33246.18

# WVR row 113
$ws8.Range("H113").Value = 1986
$ws8.Range("I113").Value = 1499.5
$ws8.Range("J113").Value = 2634.6667
$ws8.Range("K113").Value = 4498.5
$ws8.Range("L113").Value = 7904.000100000001
$ws8.Range("M113").Value = -2328.5
$ws8.Range("N113").Value = -12244.0001

# WVR row 122
$ws8.Range("H122").Value = 156251660
$ws8.Range("I122").Value = 156251660
$ws8.Range("J122").Value = 0
$ws8.Range("K122").Value = 468754980
$ws8.Range("L122").Value = 0
$ws8.Range("M122").Value = -468752530
$ws8.Range("N122").ClearContents()

# WVR row 132
$ws8.Range("H132").Value = 2942.04
$ws8.Range("I132").Value = 2877.6875
$ws8.Range("J132").Value = 3056.4443
$ws8.Range("K132").Value = 8633.0625
$ws8.Range("L132").Value = 9169.332900000001
$ws8.Range("M132").Value = -6103.0625
$ws8.Range("N132").Value = -14229.3329
